# Weekly update: insert two new price rows into the Jengibre (ginger) price
# history table on the active sheet. The table is sorted with the newest
# observation first, so the new rows are inserted near the top and roughly
# 4/5ths of the way down the table (matching the source diff), pushing the
# existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at position 241 -----------------------------------
$ws.Rows.Item(241).Insert()

$ws.Range("A241").Value = 10
$ws.Range("B241").Value = "Vega Modelo de Temuco"
$ws.Range("C241").Value = "La Araucanía"
$ws.Range("D241").Value = 45120
$ws.Range("E241").Value = 9
$ws.Range("F241").Value = 100114007
$ws.Range("G241").Value = "Jengibre"
$ws.Range("H241").Value = "Sin especificar"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 100
$ws.Range("K241").Value = 24000
$ws.Range("L241").Value = 24000
$ws.Range("M241").Value = 24000
$ws.Range("N241").Value = "$/caja 13 kilos"
$ws.Range("O241").Value = "Perú"
$ws.Range("P241").Value = 1846
$ws.Range("Q241").Value = 13
$ws.Range("R241").Value = "Hortaliza"

# --- Insert a second new row at (the now-shifted) position 322 --------
$ws.Rows.Item(322).Insert()

$ws.Range("A322").Value = 10
$ws.Range("B322").Value = "Vega Modelo de Temuco"
$ws.Range("C322").Value = "La Araucanía"
$ws.Range("D322").Value = 45121
$ws.Range("E322").Value = 9
$ws.Range("F322").Value = 100114007
$ws.Range("G322").Value = "Jengibre"
$ws.Range("H322").Value = "Sin especificar"
$ws.Range("I322").Value = "Primera"
$ws.Range("J322").Value = 80
$ws.Range("K322").Value = 24000
$ws.Range("L322").Value = 24000
$ws.Range("M322").Value = 24000
$ws.Range("N322").Value = "$/caja 13 kilos"
$ws.Range("O322").Value = "Perú"
$ws.Range("P322").Value = 1846
$ws.Range("Q322").Value = 13
$ws.Range("R322").Value = "Hortaliza"
